# Apply edits described by the diff:
# 1. Fix a few time-format typos in column C (rows 8, 9, 11, 12, 24->25)
# 2. Insert a new row at row 22 ("Narre Warren" / Tamarind 8) which pushes
#    the existing rows 22-26 down to 23-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Minor text corrections (before the row insert, so row numbers match) ---
$ws.Range("C8").Value = "31/12/20 3:30pm-5:30pm"
$ws.Range("C9").Value = "29/12/20 9:30am-10:45am"
$ws.Range("C11").Value = "29/12/20 5:30pm-05:50pm"
$ws.Range("C12").Value = "29/12/20 3:30pm-04:30pm"

# Correct the Sandringham Line exposure-period typo (currently on row 24,
# "28/12/20 7pm -7.50pm" -> "28/12/20 7:00pm-7.50pm")
$ws.Range("C24").Value = "28/12/20 7:00pm-7.50pm"

# --- Insert a new row before row 22 and populate it ---
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = "Narre Warren"
$ws.Range("B22").Value = "Tamarind 8  7b/420 Princes Highway, Narre Warren"
$ws.Range("C22").Value = "30/12/20, 6:30pm-7:00pm"
$ws.Range("D22").Value = "Case got take away"
